$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to match Modflow dictionaries (df_flowline / df_particle)
$ws.Range("I1").Value = "redox"
$ws.Range("K1").Value = "travel_distance"
$ws.Range("L1").Value = "porosity"
$ws.Range("P1").Value = "solid_density"

# Update K column (travel_distance) values for each target_aquifer row
# to equal the travel distance (xcoord) accumulated for that flowline
$ws.Range("K5").Value = 5.45045333054815
$ws.Range("K9").Value = 17.23584680498276
$ws.Range("K13").Value = 54.5045333054815
$ws.Range("K17").Value = 121.875841552958
$ws.Range("K21").Value = 172.3584680498276
$ws.Range("K25").Value = 243.7516831059159
$ws.Range("K29").Value = 298.5336237770383
$ws.Range("K33").Value = 344.7169360996551
$ws.Range("K37").Value = 385.40525105714
$ws.Range("K41").Value = 422.1902995698748
$ws.Range("K45").Value = 456.0176428159155
$ws.Range("K49").Value = 487.5033662118318
$ws.Range("K53").Value = 517.0754041494827
$ws.Range("K57").Value = 545.045333054815
$ws.Range("K61").Value = 571.6483679617371
$ws.Range("K65").Value = 597.0672475540766
$ws.Range("K69").Value = 621.447294314075
$ws.Range("K73").Value = 644.9063351516777
$ws.Range("K77").Value = 667.5414763348052
$ws.Range("K81").Value = 689.4338721993103
$ws.Range("K85").Value = 710.6521692390858
$ws.Range("K89").Value = 731.2550493177478
$ws.Range("K93").Value = 751.2931442926822
$ws.Range("K97").Value = 770.81050211428
$ws.Range("K101").Value = 789.8457265049625
$ws.Range("K105").Value = 808.4328748799342
$ws.Range("K109").Value = 826.6021743832512
$ws.Range("K113").Value = 844.3805991397495
$ws.Range("K117").Value = 861.7923402491377
$ws.Range("K121").Value = 878.8591919190293
$ws.Range("K125").Value = 895.6008713311152
$ws.Range("K129").Value = 912.0352856318312
$ws.Range("K133").Value = 928.1787563535482
$ws.Range("K137").Value = 944.04620927924
$ws.Range("K141").Value = 959.6513360398068
$ws.Range("K145").Value = 975.0067324236636
$ws.Range("K149").Value = 990.1240173735575
$ws.Range("K153").Value = 1005.013935867775
$ws.Range("K157").Value = 1019.686448275606
$ws.Range("K161").Value = 1034.150808298965
$ws.Range("K165").Value = 1048.415631233094
$ws.Range("K169").Value = 1062.488953976638
$ws.Range("K173").Value = 1076.378287978187
$ws.Range("K177").Value = 1090.09066610963
$ws.Range("K181").Value = 1103.632684296646
$ws.Range("K185").Value = 1117.010538605748
$ws.Range("K189").Value = 1130.230058379615
$ws.Range("K193").Value = 1143.296735923474
$ws.Range("K197").Value = 1156.21575317142
$ws.Range("K201").Value = 1168.992005699884
$ws.Range("K205").Value = 1181.630124403877
$ws.Range("K209").Value = 1194.134495108153
$ws.Range("K213").Value = 1206.509276348793
$ws.Range("K217").Value = 1218.75841552958
$ws.Range("K221").Value = 1230.885663631133
$ws.Range("K225").Value = 1242.89458862815
$ws.Range("K229").Value = 1254.78858775079
$ws.Range("K233").Value = 1266.570898709624
$ws.Range("K237").Value = 1278.244609989257
$ws.Range("K241").Value = 1289.812670303355
$ws.Range("K245").Value = 1301.277897293101
$ws.Range("K249").Value = 1312.64298554178
$ws.Range("K253").Value = 1323.91051397008
$ws.Range("K257").Value = 1335.08295266961
$ws.Range("K261").Value = 1346.162669225933
$ws.Range("K265").Value = 1357.151934576955
$ws.Range("K269").Value = 1368.052928447747
$ws.Range("K273").Value = 1378.867744398621
$ws.Range("K277").Value = 1389.59839451959
$ws.Range("K281").Value = 1400.246813801019
$ws.Range("K285").Value = 1410.814864207345
$ws.Range("K289").Value = 1421.304338478172
$ws.Range("K293").Value = 1431.7169636787
$ws.Range("K297").Value = 1442.054404519413
$ws.Range("K301").Value = 1452.318266463079
$ws.Range("K305").Value = 1462.510098635495
$ws.Range("K309").Value = 1472.631396554919
$ws.Range("K313").Value = 1482.683604693791
$ws.Range("K317").Value = 1492.668118885192
$ws.Range("K321").Value = 1502.586288585364
$ws.Range("K325").Value = 1512.4394190027
$ws.Range("K329").Value = 1522.228773102685
$ws.Range("K333").Value = 1531.955573497544
$ws.Range("K337").Value = 1541.62100422856
$ws.Range("K341").Value = 1551.226212448448
$ws.Range("K345").Value = 1560.772310010541
$ws.Range("K349").Value = 1570.260374971019
$ws.Range("K353").Value = 1579.691453009925
$ws.Range("K357").Value = 1589.066558776281
$ws.Range("K361").Value = 1598.386677162186
$ws.Range("K365").Value = 1607.652764510439
$ws.Range("K369").Value = 1616.865749759868
$ws.Range("K373").Value = 1626.026535532254
$ws.Range("K377").Value = 1635.135999164445
$ws.Range("K381").Value = 1644.194993689006
$ws.Range("K385").Value = 1653.204348766502
$ws.Range("K389").Value = 1662.164871572299
$ws.Range("K393").Value = 1671.07734764057
$ws.Range("K397").Value = 1679.942541667996
$ws.Range("K401").Value = 1688.761198279499
$ws.Range("K405").Value = 1697.534042758169
$ws.Range("K409").Value = 1706.261781741411
$ws.Range("K413").Value = 1714.945103885211
$ws.Range("K417").Value = 1719.27031908717
$ws.Range("K421").Value = 1723.498499109662
$ws.Range("K425").Value = 1985.659025940517
$ws.Range("K429").Value = 2247.819552771371
$ws.Range("K433").Value = 2509.980079602226
